# TC_144.xlsx edit: rename sheet, tweak a couple of labels/values, and
# widen the custom number format used by the downloaded data block.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Rename the worksheet tab from "My Series" to "Data".
$ws.Name = "Data"

# 2) F1 header label: drop "SAR (China)" from the series name.
$ws.Range("F1").Value = "(DC)Hong Kong Retail Bonds: Price: Mid: HK Link A: 07-05-2009: 3.60%"

# 3) A11 label rename.
$ws.Range("A11").Value = "Function Information"

# 4) F14 "Last Update Time" value corrected to match the rest of the row.
$ws.Cells.Item(14, 6).Value = 41781

# 5) Tiny floating point corrections on the Skewness/Kurtosis rows.
$ws.Cells.Item(20, 3).Value = 0.0156159028969747
$ws.Cells.Item(20, 4).Value = 0.5088159326089575
$ws.Cells.Item(21, 2).Value = -0.4005342754007675
$ws.Cells.Item(21, 3).Value = -1.075751539480832
$ws.Cells.Item(21, 4).Value = -0.497756315366257

# 6) Widen the custom numeric format (numFmtId 166) used by the downloaded
# data block (B27:F2695) from "0.000" to "###0.000".
$ws.Range("B27:F2695").NumberFormat = "###0.000"
